$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.757.79'
$ws.Range("E2").Value = '  -2.44%  '
$ws.Range("D3").Value = '1.565.36'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''206.33'
$ws.Range("E5").Value = '  -0.81%  '
$ws.Range("E6").Value = '  -2.17%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("E9").Value = '  -0.46%  '
$ws.Range("D10").Value = '''0.0584'
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = '1.787.23'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '1.561.20'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = '''61.58'
$ws.Range("E16").Value = '  -2.87%  '
$ws.Range("D17").Value = '26.786.29'
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").Value = '''214.22'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").Value = '''7.35'
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '''4.09'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").Value = '''9.32'
$ws.Range("E23").Value = '  -2.04%  '
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").Value = '''152.59'
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("D27").Value = '''14.91'
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("E30").Value = '  -3.48%  '
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("E32").Value = '  -1.40%  '
$ws.Range("D33").Value = '1.383.43'
$ws.Range("E33").Value = '  +1.02%  '
$ws.Range("D34").Value = '''2.93'
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("E35").Value = '  +1.11%  '
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("D37").Value = '''0.923'
$ws.Range("E37").Value = '  -3.36%  '
$ws.Range("E38").Value = '  -2.39%  '
$ws.Range("D39").Value = '''0.526'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '''0.817'
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").Value = '''5.34'
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("D46").Value = '''63.27'
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("D47").Value = '1.700.45'
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").Value = '''85.35'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").Value = '0.0₇0990'
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").Value = '''0.0951'
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("E51").Value = '  -0.62%  '
